$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 28: was the "new" kitchen trial, becomes the "catch" trial ---
# Remove H28 ("kitchens") entirely - row 28 no longer has a category value
$ws.Range("H28").ClearContents()

# J28: "new" -> "catch"
$ws.Range("J28").Value = "catch"

# L28: stimulus image path updated
$ws.Range("L28").Value = "stimuli/catch_16.jpg"

# Remove the now-irrelevant numeric stat columns M28:V28
$ws.Range("M28:V28").ClearContents()

# --- Row 29: was the "catch" trial, becomes the "new" kitchen trial ---
# Add H29 ("kitchens") category value
$ws.Range("H29").Value = "kitchens"

# J29: "catch" -> "new"
$ws.Range("J29").Value = "new"

# L29: stimulus image path updated
$ws.Range("L29").Value = "stimuli/img_s9are.png"

# Populate the numeric stat columns M29:V29 (moved down from former row 28)
$ws.Range("M29").Value = 90.14285714285714
$ws.Range("N29").Value = 75.22857142857143
$ws.Range("O29").Value = 82.68571428571428
$ws.Range("P29").Value = 35
$ws.Range("Q29").Value = 10
$ws.Range("R29").Value = 10
$ws.Range("S29").Value = 10
$ws.Range("T29").Value = 10
$ws.Range("U29").Value = 10
$ws.Range("V29").Value = 10
